$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 6.2920276419332239
$ws.Range("A3").Value = 6.0490642775245593
$ws.Range("A4").Value = 4.7342902645750051
$ws.Range("B4").Value = 4.9000000000000004
$ws.Range("C4").Value = 4.871999999999999
$ws.Range("D4").Value = 2.4500000000000002
$ws.Range("A5").Value = 4.0146724957974183
$ws.Range("B5").Value = 4.2
$ws.Range("C5").Value = 4.0508695652173898
$ws.Range("D5").Value = 2.66
$ws.Range("A6").Value = 3.9453446815050182
$ws.Range("B6").Value = 3.64
$ws.Range("C6").Value = 4.0981818181818177
$ws.Range("D6").Value = 2.87
$ws.Range("A7").Value = 3.8198181722706677
$ws.Range("B7").Value = 3.43
$ws.Range("C7").Value = 3.01
$ws.Range("D7").Value = 3.01
$ws.Range("E7").Value = 3.01
$ws.Range("A8").Value = 3.5034182301568149
$ws.Range("B8").Value = 3.29
$ws.Range("C8").Value = 4.0981818181818177
$ws.Range("D8").Value = 2.87
$ws.Range("A9").Value = 2.7812219848714714
$ws.Range("B9").Value = 2.87
$ws.Range("C9").Value = 4.0508695652173898
$ws.Range("D9").Value = 2.66
$ws.Range("A10").Value = 1.9911650048790033
$ws.Range("B10").Value = 1.4
$ws.Range("C10").Value = 1.3248148148148147
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 3.22
$ws.Range("A11").Value = 1.6128359075221363
$ws.Range("B11").Value = 1.1200000000000001
$ws.Range("C11").Value = 1.02
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 3.08
$ws.Range("A12").Value = 1.5026277781444826
$ws.Range("B12").Value = 1.05
$ws.Range("C12").Value = 0.64749999999999996
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 3.01
$ws.Range("A13").Value = 3.9918979353713615
$ws.Range("B13").Value = 3.99
$ws.Range("C13").Value = 3.99
$ws.Range("D13").Value = 3.99
$ws.Range("E13").Value = 3.99
$ws.Range("A14").Value = 3.9915858906130106
$ws.Range("B14").Value = 3.99
$ws.Range("C14").Value = 3.99
$ws.Range("D14").Value = 3.85
$ws.Range("E14").Value = 4.13
$ws.Range("A15").Value = 3.9942603905120166
$ws.Range("B15").Value = 3.99
$ws.Range("C15").Value = 3.99
$ws.Range("D15").Value = 3.64
$ws.Range("E15").Value = 4.34
$ws.Range("A16").Value = 4.2899705318464783
$ws.Range("B16").Value = 4.2
$ws.Range("C16").Value = 3.9900000000000007
$ws.Range("D16").Value = 3.43
$ws.Range("E16").Value = 4.55
$ws.Range("A17").Value = 4.9123790596025163
$ws.Range("A18").Value = 4.9863967766256838
$ws.Range("A19").Value = 5.0546361570371774
$ws.Range("A20").Value = 5.7723631572360761
$ws.Range("A21").Value = 6.1787883259107117
$ws.Range("A22").Value = 6.2609135440857155
$ws.Range("A23").Value = 3.3331391742079211
$ws.Range("B23").Value = 3.36
$ws.Range("C23").Value = 3.5
$ws.Range("D23").Value = 3.01
$ws.Range("A24").Value = 3.5067150619471064
$ws.Range("B24").Value = 3.5
$ws.Range("C24").Value = 3.4999999999999996
$ws.Range("D24").Value = 2.4500000000000002
$ws.Range("E24").Value = 4.55
$ws.Range("A25").Value = 3.374803063316977
$ws.Range("B25").Value = 3.29
$ws.Range("C25").Value = 3.0100000000000002
$ws.Range("D25").Value = 2.87
$ws.Range("E25").Value = 3.15
$ws.Range("A26").Value = 2.9455104173790643
$ws.Range("B26").Value = 3.08
$ws.Range("C26").Value = 3.0100000000000002
$ws.Range("D26").Value = 2.87
$ws.Range("E26").Value = 3.15
$ws.Range("A27").Value = 1.8951320462239079
$ws.Range("B27").Value = 1.26
$ws.Range("C27").Value = 0.62999999999999989
$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 1.26
